$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.323.72'
$ws.Range("E2").Value = '  +0.15%  '

$ws.Range("D3").Value = '3.483.11'
$ws.Range("E3").Value = '  -0.89%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '596.94'
$ws.Range("E5").Value = '  +0.23%  '

$ws.Range("D6").Value = '177.17'
$ws.Range("E6").Value = '  +2.49%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  -0.25%  '

$ws.Range("E9").Value = '  -2.18%  '

$ws.Range("D10").Value = '7.08'
$ws.Range("E10").Value = '  -2.94%  '

$ws.Range("D11").Value = '0.425'
$ws.Range("E11").Value = '  -2.45%  '

$ws.Range("D12").Value = '4.083.49'
$ws.Range("E12").Value = '  -1.09%  '

$ws.Range("D13").Value = '31.73'
$ws.Range("E13").Value = '  +10.47%  '

$ws.Range("E14").Value = '  -0.18%  '

$ws.Range("D15").Value = '67.258.55'
$ws.Range("E15").Value = '  +0.09%  '

$ws.Range("D17").Value = '3.487.08'
$ws.Range("E17").Value = '  -0.82%  '

$ws.Range("E18").Value = '  -1.43%  '

$ws.Range("D19").Value = '14.41'
$ws.Range("E19").Value = '  +1.16%  '

$ws.Range("D20").Value = '388.86'
$ws.Range("E20").Value = '  -1.88%  '

$ws.Range("D21").Value = '7.92'
$ws.Range("E21").Value = '  -0.79%  '

$ws.Range("D22").Value = '72.86'
$ws.Range("E22").Value = '  -0.71%  '

$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("E24").Value = '  -0.87%  '

$ws.Range("D25").Value = '5.71'
$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("D26").Value = '0.0000121'
$ws.Range("E26").Value = '  -0.68%  '

$ws.Range("D27").Value = '10.26'
$ws.Range("E27").Value = '  -0.44%  '

$ws.Range("D28").Value = '0.177'
$ws.Range("E28").Value = '  -1.67%  '

$ws.Range("D29").Value = '0.996'
$ws.Range("E29").Value = '  -0.31%  '

$ws.Range("D30").Value = '6.17'
$ws.Range("E30").Value = '  -1.81%  '

$ws.Range("D31").Value = '1.42'
$ws.Range("E31").Value = '  -2.57%  '

$ws.Range("E32").Value = '  -1.80%  '

$ws.Range("D33").Value = '23.56'
$ws.Range("E33").Value = '  -2.25%  '

$ws.Range("D34").Value = '7.27'
$ws.Range("E34").Value = '  -1.81%  '

$ws.Range("E35").Value = '  -0.26%  '

$ws.Range("D36").Value = '163.68'
$ws.Range("E36").Value = '  -0.11%  '

$ws.Range("E37").Value = '  +0.67%  '

$ws.Range("D38").Value = '0.869'
$ws.Range("E38").Value = '  -2.90%  '

$ws.Range("D39").Value = '6.98'
$ws.Range("E39").Value = '  +0.93%  '

$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D40").Value = '27.28'
$ws.Range("E40").Value = '  -0.42%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '4.63'
$ws.Range("E41").Value = '  -1.96%  '

$ws.Range("D42").Value = '26.34'
$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("D43").Value = '2.812.70'
$ws.Range("E43").Value = '  -0.11%  '

$ws.Range("D44").Value = '0.0721'
$ws.Range("E44").Value = '  -3.44%  '

$ws.Range("D45").Value = '2.57'
$ws.Range("E45").Value = '  -1.93%  '

$ws.Range("D46").Value = '42.22'

$ws.Range("D47").Value = '341.83'
$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("E48").Value = '  -3.41%  '

$ws.Range("E49").Value = '  -2.92%  '

$ws.Range("D50").Value = '33.31'
$ws.Range("E50").Value = '  -2.14%  '

$ws.Range("E51").Value = '  -2.29%  '
